$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7751
$ws1.Range("F4").Value = 7914
$ws1.Range("F8").Value = 6724
$ws1.Range("F9").Value = 6724
$ws1.Range("F13").Value = 45
$ws1.Range("F18").Value = 474
$ws1.Range("F24").Value = 3883
$ws1.Range("F26").Value = 374
$ws1.Range("F28").Value = 288
$ws1.Range("F29").Value = 1502
$ws1.Range("F31").Value = 63
$ws1.Range("F32").Value = 2776
$ws1.Range("F33").Value = 1899
$ws1.Range("F34").Value = 38
$ws1.Range("F37").Value = 59
$ws1.Range("F38").Value = 3746
$ws1.Range("F39").Value = 335
$ws1.Range("F42").Value = 925
$ws1.Range("F43").Value = 552
$ws1.Range("F48").Value = 562
$ws1.Range("F50").Value = 10

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 13
$ws2.Range("F17").Value = 140

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 140

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 140
$ws4.Range("F7").Value = 7751
$ws4.Range("F9").Value = 7914
$ws4.Range("F12").Value = 6724
$ws4.Range("F16").Value = 45
$ws4.Range("F21").Value = 474
$ws4.Range("F26").Value = 3883
$ws4.Range("F30").Value = 374
$ws4.Range("F32").Value = 1502
$ws4.Range("F34").Value = 63
$ws4.Range("F35").Value = 2776
$ws4.Range("F36").Value = 1899
$ws4.Range("F37").Value = 38
$ws4.Range("F41").Value = 335
$ws4.Range("F45").Value = 925
$ws4.Range("F46").Value = 552
$ws4.Range("F47").Value = 140
$ws4.Range("F49").Value = 562
